$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("schedule")

# Add the new note about Edward covering the 10am class to the existing
# "data viz lab" and "I won't be here" comments (rows 15 & 16, column F).
$ws.Range("F15").Value = "data viz lab`nEdward cover 10am class"
$ws.Range("F16").Value = "I won't be here`nEdward cover 10am class"

# Insert the new holiday into the topic column (D), pushing the existing
# topics for "Logistic Regression" and "Analysis & Poster work" down one
# week, and note that the prior week's topic could run long.
$ws.Range("E36").Value = "could take an entire class period"
$ws.Range("D39").Value = $ws.Range("D38").Value()
$ws.Range("D38").Value = $ws.Range("D37").Value()
$ws.Range("D37").Value = "No Class - Holiday"

# Update the frozen-pane scroll position and active selection to match
# where the editor was working.
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("F17").Select()
